# "Generate Report for Handoff"
# The two tracked e2e markdown files were just handed off for a new round of
# localization: a brand new source file (e7124884-...) replaces the file that
# used to be 6a9f4aba-..., and a second file rolls over to a fresh GUID
# (ffff65dcf7da-...) in place of ea05824d-.... Their status flips from
# "Handed back: in sync with en-US" to "Ready for handoff", new handoff
# xliffs/timestamps are recorded, and the (now stale) handback info is
# cleared back out to blank until the new round comes back.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# ---- new values -----------------------------------------------------------
$newFile1        = "e7124884-0877-4844-a235-bea08737dd45.md"
$newPath1        = "e2e\e7124884-0877-4844-a235-bea08737dd45.md"
$newFile2        = "ffff65dcf7da-4ad3-4385-a7c7-8686eed17529.md"
$newPath2        = "e2e\ffff65dcf7da-4ad3-4385-a7c7-8686eed17529.md"

$newStatus       = "Ready for handoff"
$newHoDate       = "2016-09-07 01:22:41"

$newHandoffZh    = "e7124884-0877-4844-a235-bea08737dd45.1e2602b1917371dd72aa01bc3efb50038639f6c5.zh-cn.xlf"
$newHandoffDe    = "e7124884-0877-4844-a235-bea08737dd45.1e2602b1917371dd72aa01bc3efb50038639f6c5.de-de.xlf"
$newHandoffDtZh  = "2016-09-07 01:22:36"
$newHandbackDt   = "0001-01-01 00:00:00"

# =============================================================================
# Overview sheet
# =============================================================================
$wsOverview.Range("A2").Value = $newFile1
$wsOverview.Range("B2").Value = $newPath1
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("G2").Value = $newHoDate

$wsOverview.Range("A3").Value = $newFile2
$wsOverview.Range("B3").Value = $newPath2
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsOverview.Range("G3").Value = $newHoDate

foreach ($h in @($wsOverview.Hyperlinks)) {
    $addr = $h.Range.Address()
    if ($addr -eq '$B$2') {
        $h.TextToDisplay = $newPath1
    } elseif ($addr -eq '$B$3') {
        $h.TextToDisplay = $newPath2
    }
}

$wsOverview.Columns.Item(5).ColumnWidth = 17.22
$wsOverview.Columns.Item(6).ColumnWidth = 17.22

# =============================================================================
# zh-cn sheet
# =============================================================================
$wsZhCn.Range("A2").Value = $newFile1
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("G2").Value = $newHandoffZh
$wsZhCn.Range("H2").Value = $newHandoffDtZh
$wsZhCn.Range("I2").Value = ""
$wsZhCn.Range("J2").Value = ""
$wsZhCn.Range("K2").Value = $newHandbackDt

$wsZhCn.Range("A3").Value = $newFile2
$wsZhCn.Range("C3").Value = $newStatus
$wsZhCn.Range("F3").Value = "True"
$wsZhCn.Range("G3").Value = $newHandoffZh
$wsZhCn.Range("H3").Value = $newHandoffDtZh
$wsZhCn.Range("I3").Value = ""
$wsZhCn.Range("J3").Value = ""
$wsZhCn.Range("K3").Value = $newHandbackDt

foreach ($h in @($wsZhCn.Hyperlinks)) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') {
        $h.TextToDisplay = $newFile1
    } elseif ($addr -eq '$A$3') {
        $h.TextToDisplay = $newFile2
    }
}
foreach ($h in @($wsZhCn.Hyperlinks)) {
    $addr = $h.Range.Address()
    if ($addr -eq '$I$2' -or $addr -eq '$I$3') {
        $h.Delete()
    }
}

$wsZhCn.Range("I2").Style = "Normal"
$wsZhCn.Range("J2").Style = "Normal"
$wsZhCn.Range("I3").Style = "Normal"
$wsZhCn.Range("J3").Style = "Normal"

$wsZhCn.Columns.Item(3).ColumnWidth = 17.22
$wsZhCn.Columns.Item(9).ColumnWidth = 18.67
$wsZhCn.Columns.Item(10).ColumnWidth = 21.71

# =============================================================================
# de-de sheet
# =============================================================================
$wsDeDe.Range("A2").Value = $newFile1
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("G2").Value = $newHandoffDe
$wsDeDe.Range("H2").Value = $newHoDate
$wsDeDe.Range("I2").Value = ""
$wsDeDe.Range("J2").Value = ""
$wsDeDe.Range("K2").Value = $newHandbackDt

$wsDeDe.Range("A3").Value = $newFile2
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Range("F3").Value = "True"
$wsDeDe.Range("G3").Value = $newHandoffDe
$wsDeDe.Range("H3").Value = $newHoDate
$wsDeDe.Range("I3").Value = ""
$wsDeDe.Range("J3").Value = ""
$wsDeDe.Range("K3").Value = $newHandbackDt

foreach ($h in @($wsDeDe.Hyperlinks)) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') {
        $h.TextToDisplay = $newFile1
    } elseif ($addr -eq '$A$3') {
        $h.TextToDisplay = $newFile2
    }
}
foreach ($h in @($wsDeDe.Hyperlinks)) {
    $addr = $h.Range.Address()
    if ($addr -eq '$I$2' -or $addr -eq '$I$3') {
        $h.Delete()
    }
}

$wsDeDe.Range("I2").Style = "Normal"
$wsDeDe.Range("J2").Style = "Normal"
$wsDeDe.Range("I3").Style = "Normal"
$wsDeDe.Range("J3").Style = "Normal"

$wsDeDe.Columns.Item(3).ColumnWidth = 17.22
$wsDeDe.Columns.Item(9).ColumnWidth = 18.67
$wsDeDe.Columns.Item(10).ColumnWidth = 21.71

Write-Host "Applied handoff report update."
